# ---------------------------------------------------------------------------
# Add files via upload
#
# 1. Zero-out every "PRESUPUESTO" (budget) figure on the "VENTA MENSUAL"
#    sheet (column G, rows 2-55, including the total in row 55).
# 2. Add a brand-new "CUMPLIMIENTO MENSUAL" worksheet (after "VENTA MENSUAL")
#    that breaks the same budget/sale numbers down by GRUPO instead of by
#    CLIENTE, with a POR CUMPLIR (remaining) and CUMPLIMIENTO (% achieved)
#    column, plus a TOTAL row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: VENTA MENSUAL -> column G (PRESUPUESTO) becomes 0 everywhere.
# ---------------------------------------------------------------------------
$ventaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
for ($r = 2; $r -le 55; $r++) {
    $ventaMensual.Cells.Item($r, 7).Value = 0
}

# ---------------------------------------------------------------------------
# Step 2: new "CUMPLIMIENTO MENSUAL" sheet, placed right after "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Column widths (character units) -- tuned so the stored OOXML width lands on
# 26 / 22 / 17 / 13 / 17 / 26, matching the source column layout.
$ws.Columns.Item(1).ColumnWidth = 25.1
$ws.Columns.Item(2).ColumnWidth = 21.1
$ws.Columns.Item(3).ColumnWidth = 16.1
$ws.Columns.Item(4).ColumnWidth = 12.1
$ws.Columns.Item(5).ColumnWidth = 16.1
$ws.Columns.Item(6).ColumnWidth = 25.1

# Header row
$headers = @("ASESOR", "GRUPO", "PRESUPUESTO", "VENTA", "POR CUMPLIR", "CUMPLIMIENTO")
for ($c = 1; $c -le 6; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$asesor = "CASTRO ALCIVAR EDA MARIA"

# GRUPO, PRESUPUESTO, VENTA  (POR CUMPLIR = PRESUPUESTO-VENTA, CUMPLIMIENTO = VENTA/PRESUPUESTO)
$grupos = @(
    @("240X120 PORCELANATO", 5820, 1669.25),
    @("240X80 PORCELANATO", 13728, 0),
    @("FREGADEROS DE COCINA", 646, 128.74),
    @("GRANITO", 238.32, 0),
    @("GRIFERIAS", 106.82, 0),
    @("INODOROS", 2100, 0),
    @("LAVABOS", 1000, 0),
    @("LED", 300, 0),
    @("NO RESURTIBLES", 1300.5, 0),
    @("OTROS", 0, 0),
    @("PANELES DECORATIVOS", 350, 0),
    @("PANELES PU", 230, 0),
    @("PANELES PVC", 966, 0),
    @("PIEDRA SINTERIZADA", 15690, 2568.3),
    @("PORCELANATO", 45745.689, 346.47),
    @("PUERTAS DE SEGURIDAD", 1142, 0),
    @("SAL SOLUBLE", 1600, 0)
)

$row = 2
$totalPresupuesto = 0
$totalVenta = 0
foreach ($g in $grupos) {
    $nombre = $g[0]
    $presupuesto = $g[1]
    $venta = $g[2]
    $porCumplir = $presupuesto - $venta

    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $nombre

    $cC = $ws.Cells.Item($row, 3)
    $cC.Value = $presupuesto
    $cC.NumberFormat = "`"$`"#,##0.00"

    $cD = $ws.Cells.Item($row, 4)
    $cD.Value = $venta
    $cD.NumberFormat = "`"$`"#,##0.00"

    $cE = $ws.Cells.Item($row, 5)
    $cE.Value = $porCumplir
    $cE.NumberFormat = "`"$`"#,##0.00"

    $cF = $ws.Cells.Item($row, 6)
    if ($presupuesto -ne 0) {
        $cF.Value = $venta / $presupuesto
    } else {
        $cF.Value = 0
    }
    $cF.NumberFormat = "0.00%"

    $totalPresupuesto += $presupuesto
    $totalVenta += $venta

    $row++
}

# TOTAL row
$totalPorCumplir = $totalPresupuesto - $totalVenta
$ws.Cells.Item($row, 2).Value = "TOTAL"
$ws.Cells.Item($row, 2).HorizontalAlignment = -4152

$cC = $ws.Cells.Item($row, 3)
$cC.Value = $totalPresupuesto
$cC.NumberFormat = "`"$`"#,##0.00"

$cD = $ws.Cells.Item($row, 4)
$cD.Value = $totalVenta
$cD.NumberFormat = "`"$`"#,##0.00"

$cE = $ws.Cells.Item($row, 5)
$cE.Value = $totalPorCumplir
$cE.NumberFormat = "`"$`"#,##0.00"

$cF = $ws.Cells.Item($row, 6)
$cF.Value = $totalVenta / $totalPresupuesto
$cF.NumberFormat = "0.00%"

$ws.Range("A1").Select()
